$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data correction: End Location for row 2 (Atm_on) shrinks by 1bp ---
$ws.Range("C2").Value = 53504218

# --- New header cell I1: "Primer_Ident" ---
$ws.Range("I1").Value = "Primer_Ident"

# --- New Primer_Ident numeric values for each data row (default 90, two exceptions at 100) ---
for ($r = 2; $r -le 41; $r++) {
    $val = 90
    if ($r -eq 7 -or $r -eq 16) { $val = 100 }
    $ws.Cells.Item($r, 9).Value = $val
}

# --- Match formatting: header I1/J1 take on the same style as the other header cells (e.g. B1) ---
$ws.Range("B1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# --- Match formatting: data rows I/J take on the same style as column B's "primary" style (s=4) ---
$ws.Range("B2").Copy()
$ws.Range("I2:I41").PasteSpecial(-4122)
$ws.Range("J2:J41").PasteSpecial(-4122)

# --- Row heights: all rows (header + data) shrink slightly from 15.75 to 15 ---
for ($r = 1; $r -le 41; $r++) {
    $ws.Rows.Item($r).RowHeight = 15
}

# --- Leave selection where the author left off ---
$ws.Range("M28").Select()
